$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (styles) of the last existing data row (94) down into
# the two new rows (95, 96) so column A keeps its bold/border style and
# column E keeps its date number format.
$ws.Range("A94:V94").Copy($ws.Range("A95:V95"))
$ws.Range("A94:V94").Copy($ws.Range("A96:V96"))

# ---- Row 95 ----
# Columns B (paraguay), C (primera-division) and D (2023) already have the
# correct text values after the row copy above, so they are left untouched
# to preserve their inline-string cell type.
$ws.Cells.Item(95,1).Value  = 94
$ws.Cells.Item(95,5).Value  = 45223.95833333334
$ws.Cells.Item(95,6).Value  = "Sportivo Trinidense"
$ws.Cells.Item(95,7).Value  = 2
$ws.Cells.Item(95,8).Value  = "Cerro Porteno"
$ws.Cells.Item(95,9).Value  = 5
$ws.Cells.Item(95,10).Value = 3.87
$ws.Cells.Item(95,11).Value = "22/10/2023 00:42"
$ws.Cells.Item(95,12).Value = 4.12
$ws.Cells.Item(95,13).Value = "24/10/2023 22:55"
$ws.Cells.Item(95,14).Value = 3.65
$ws.Cells.Item(95,15).Value = "22/10/2023 00:42"
$ws.Cells.Item(95,16).Value = 3.59
$ws.Cells.Item(95,17).Value = "24/10/2023 22:50"
$ws.Cells.Item(95,18).Value = 1.95
$ws.Cells.Item(95,19).Value = "22/10/2023 00:42"
$ws.Cells.Item(95,20).Value = 1.94
$ws.Cells.Item(95,21).Value = "24/10/2023 22:50"
$ws.Cells.Item(95,22).Value = "https://www.betexplorer.com/football/paraguay/primera-division/sportivo-trinidense-cerro-porteno/fNuFVZ1f/"

# ---- Row 96 ----
# Columns B, C and D already carry the correct text values from the row copy.
$ws.Cells.Item(96,1).Value  = 95
$ws.Cells.Item(96,5).Value  = 45224.0625
$ws.Cells.Item(96,6).Value  = "Guarani"
$ws.Cells.Item(96,7).Value  = 0
$ws.Cells.Item(96,8).Value  = "Olimpia Asuncion"
$ws.Cells.Item(96,9).Value  = 1
$ws.Cells.Item(96,10).Value = 3.38
$ws.Cells.Item(96,11).Value = "21/10/2023 22:12"
$ws.Cells.Item(96,12).Value = 3.77
$ws.Cells.Item(96,13).Value = "25/10/2023 00:16"
$ws.Cells.Item(96,14).Value = 3.52
$ws.Cells.Item(96,15).Value = "21/10/2023 22:12"
$ws.Cells.Item(96,16).Value = 3.5
$ws.Cells.Item(96,17).Value = "25/10/2023 01:07"
$ws.Cells.Item(96,18).Value = 2.15
$ws.Cells.Item(96,19).Value = "21/10/2023 22:12"
$ws.Cells.Item(96,20).Value = 2.07
$ws.Cells.Item(96,21).Value = "25/10/2023 00:16"
$ws.Cells.Item(96,22).Value = "https://www.betexplorer.com/football/paraguay/primera-division/guarani-olimpia-asuncion/YDvJUgH0/"
